# Update "Planilha de Atividades" - add 15/05 entry and revise 14/05 text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revise the 14/05 (row 6) activity description and taller row height ---
$ws.Range("A6").Value = "Deletei os itens do menu do sistema das páginas que foram deletada e terminei de deletar as tabelas, views e classes que ainda faltavam."
$ws.Rows("6").RowHeight = 60

# --- Add a new row 7 for 15/05, copying the formatting used by row 6 ---
$ws.Range("A6:B6").Copy()
$ws.Range("A7:B7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("A7").Value = "Cadastrei no BD novos itens no Menu de Páginas que ainda serão implementadas."
$ws.Range("B7").Value = 43235
$ws.Rows("7").RowHeight = 30

# --- Move the active selection to the next empty row, as in the source file ---
$ws.Range("A8").Select()
